$d = $word.ActiveDocument

$newText = "Vi förväntar oss att ni återkommer med ett skriftligt svar på vårt klagomål och även beskriver vilka korrigerande åtgärder ni satt in för att rätta till identifierade brister i er efterlevnad av den svenska FSC standarden."

# ---------------------------------------------------------------------------
# 1. Insert a new paragraph containing the "Vi förväntar oss..." text right
#    after the "Nedan presenteras fynd..." paragraph near the top of the
#    document.
# ---------------------------------------------------------------------------
$introIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Nedan presenteras fynd")) {
        $introIndex = $i
        break
    }
}

if ($introIndex -gt 0) {
    $introPara = $d.Paragraphs.Item($introIndex)
    $introPara.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($introIndex + 1)
    $newPara.Range.Text = $newText
}

# ---------------------------------------------------------------------------
# 2. Remove the old trailing copy of that same sentence together with the two
#    blank paragraphs that used to precede it at the very end of the
#    document (right before the final sectPr).
# ---------------------------------------------------------------------------
$trailIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Vi förväntar oss att ni återkommer")) {
        if ($i -ne ($introIndex + 1)) {
            $trailIndex = $i
        }
    }
}

if ($trailIndex -gt 0) {
    $endPara = $d.Paragraphs.Item($trailIndex)
    $startPara = $d.Paragraphs.Item($trailIndex - 2)
    $range = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $range.Delete()
}

# ---------------------------------------------------------------------------
# 3. Update the date shown in the first-page header from 2023-11-13 to
#    2023-11-14.
# ---------------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(2)
$hdr.Range.Find.Execute("2023-11-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-11-14", 2)
